# Scheduled-runner refresh: rewrites the profit-calc columns (H:N) of each
# class's leve table with freshly pulled market-board prices.
# Generated from the authoritative cell-level diff of Kujata_Profits.xlsx.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 258
$ws.Range("I33").Value = 177.66667
$ws.Range("K33").Value = 177.66667
$ws.Range("M33").Value = 51.33332999999999
# Row 40
$ws.Range("H40").Value = 2998.4285
$ws.Range("I40").Value = 2333
$ws.Range("J40").Value = 3497.5
$ws.Range("K40").Value = 2333
$ws.Range("L40").Value = 3497.5
$ws.Range("M40").Value = -2158
$ws.Range("N40").Value = -3847.5
# Row 96
$ws.Range("H96").Value = 1033.3334
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1033.3334
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3100.0002
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5846.0002
# Row 135
$ws.Range("H135").Value = 166667410
$ws.Range("I135").Value = 554.6667
$ws.Range("J135").Value = 333334270
$ws.Range("K135").Value = 4992.0003
$ws.Range("L135").Value = 3000008430
$ws.Range("M135").Value = -2457.0003
$ws.Range("N135").Value = -3000013500
# Row 137
$ws.Range("H137").Value = 1340.25
$ws.Range("I137").Value = 1307.4166
$ws.Range("J137").Value = 1438.75
$ws.Range("K137").Value = 3922.2498
$ws.Range("L137").Value = 4316.25
$ws.Range("M137").Value = -1372.2498
$ws.Range("N137").Value = -9416.25
# Row 138
$ws.Range("H138").Value = 1161.75
$ws.Range("I138").Value = 736.69446
$ws.Range("K138").Value = 2210.08338
$ws.Range("M138").Value = 2929.91662

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 1765.625
$ws.Range("I28").Value = 1765.625
$ws.Range("K28").Value = 1765.625
$ws.Range("M28").Value = -1573.625
# Row 32
$ws.Range("H32").Value = 4566.971
$ws.Range("I32").Value = 4297.879
$ws.Range("K32").Value = 4297.879
$ws.Range("M32").Value = -4010.879
# Row 61
$ws.Range("H61").Value = 1470.5238
$ws.Range("I61").Value = 1104.5
$ws.Range("K61").Value = 1104.5
$ws.Range("M61").Value = -892.5
# Row 74
$ws.Range("H74").Value = 1013.04346
$ws.Range("I74").Value = 814.3
$ws.Range("K74").Value = 814.3
$ws.Range("M74").Value = 59.70000000000005
# Row 77
$ws.Range("H77").Value = 1013.04346
$ws.Range("I77").Value = 814.3
$ws.Range("K77").Value = 4071.5
$ws.Range("M77").Value = 296.5
# Row 99
$ws.Range("H99").Value = 1765.625
$ws.Range("I99").Value = 1765.625
$ws.Range("K99").Value = 1765.625
$ws.Range("M99").Value = 1229.375
# Row 133
$ws.Range("H133").Value = 28175.809
$ws.Range("J133").Value = 28182.84
$ws.Range("L133").Value = 28182.84
$ws.Range("N133").Value = -33242.84
# Row 136
$ws.Range("H136").Value = 1470.5238
$ws.Range("I136").Value = 1104.5
$ws.Range("K136").Value = 3313.5
$ws.Range("M136").Value = -763.5

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 562.5
$ws.Range("I22").Value = 516.6667
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 516.6667
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -166.6667
$ws.Range("N22").Value = -1400
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 108
$ws.Range("H108").Value = 7200
$ws.Range("J108").Value = 7200
$ws.Range("L108").Value = 7200
$ws.Range("N108").Value = -14880

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1426.2727
$ws.Range("I5").Value = 1513.9
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 4541.700000000001
$ws.Range("L5").Value = 1650
$ws.Range("M5").Value = -4429.700000000001
$ws.Range("N5").Value = -1874
# Row 59
$ws.Range("H59").Value = 3207
$ws.Range("J59").Value = 8475
$ws.Range("L59").Value = 25425
$ws.Range("N59").Value = -26505
# Row 135
$ws.Range("H135").Value = 1426.2727
$ws.Range("I135").Value = 1513.9
$ws.Range("J135").Value = 550
$ws.Range("K135").Value = 13625.1
$ws.Range("L135").Value = 4950
$ws.Range("M135").Value = -11090.1
$ws.Range("N135").Value = -10020
# Row 138
$ws.Range("H138").Value = 2532.5
$ws.Range("I138").Value = 2532.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 7597.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -2457.5
$ws.Range("N138").ClearContents()
# Row 140
$ws.Range("H140").Value = 2213.182
$ws.Range("I140").Value = 2180.476
$ws.Range("J140").Value = 2900
$ws.Range("K140").Value = 6541.428
$ws.Range("L140").Value = 8700
$ws.Range("M140").Value = -1361.428
$ws.Range("N140").Value = -19060

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3582.6667
$ws.Range("I113").Value = 1708.8
$ws.Range("K113").Value = 1708.8
$ws.Range("M113").Value = 461.2
# Row 122
$ws.Range("H122").Value = 167834.78
$ws.Range("I122").Value = 1314.125
$ws.Range("K122").Value = 3942.375
$ws.Range("M122").Value = -1492.375
# Row 132
$ws.Range("H132").Value = 1930.579
$ws.Range("I132").Value = 1569.5883
$ws.Range("K132").Value = 4708.7649
$ws.Range("M132").Value = -2178.7649

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2324.8333
$ws.Range("J7").Value = 2242.3333
$ws.Range("L7").Value = 2242.3333
$ws.Range("N7").Value = -2466.3333
# Row 126
$ws.Range("H126").Value = 2324.8333
$ws.Range("J126").Value = 2242.3333
$ws.Range("L126").Value = 6726.999899999999
$ws.Range("N126").Value = -11666.9999
# Row 132
$ws.Range("H132").Value = 37302.605
$ws.Range("I132").Value = 1027.3529
$ws.Range("K132").Value = 3082.0587
$ws.Range("M132").Value = -552.0587000000005
# Row 136
$ws.Range("H136").Value = 1996.091
$ws.Range("I136").Value = 1779
$ws.Range("K136").Value = 5337
$ws.Range("M136").Value = -2787

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 14619.8
$ws.Range("J86").Value = 14619.8
$ws.Range("L86").Value = 14619.8
$ws.Range("N86").Value = -16865.8
# Row 89
$ws.Range("H89").Value = 14619.8
$ws.Range("J89").Value = 14619.8
$ws.Range("L89").Value = 73099
$ws.Range("N89").Value = -84331
# Row 96
$ws.Range("H96").Value = 3697
$ws.Range("I96").Value = 4040
$ws.Range("J96").Value = 2839.5
$ws.Range("K96").Value = 4040
$ws.Range("L96").Value = 2839.5
$ws.Range("M96").Value = -2667
$ws.Range("N96").Value = -5585.5
# Row 132
$ws.Range("H132").Value = 3248.8076
$ws.Range("I132").Value = 2903.8572
$ws.Range("K132").Value = 8711.571599999999
$ws.Range("M132").Value = -6181.571599999999

